# Apply numeric value updates across multiple sheets per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 255.66667
$ws.Range("I2").Value = 218.23077
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 218.23077
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -105.23077
$ws.Range("N2").Value = -725
$ws.Range("H40").Value = 6916.8335
$ws.Range("J40").Value = 9834
$ws.Range("L40").Value = 9834
$ws.Range("N40").Value = -10184
$ws.Range("H55").Value = 496.22223
$ws.Range("I55").Value = 332.83334
$ws.Range("J55").Value = 823
$ws.Range("K55").Value = 332.83334
$ws.Range("L55").Value = 823
$ws.Range("M55").Value = -118.83334
$ws.Range("N55").Value = -1251
$ws.Range("H113").Value = 8006.6665
$ws.Range("I113").Value = 25602.5
$ws.Range("J113").Value = 5299.615
$ws.Range("K113").Value = 25602.5
$ws.Range("L113").Value = 5299.615
$ws.Range("M113").Value = -22348.5
$ws.Range("N113").Value = -11807.615
$ws.Range("H116").Value = 7262
$ws.Range("I116").Value = 4913.2
$ws.Range("J116").Value = 19006
$ws.Range("K116").Value = 4913.2
$ws.Range("L116").Value = 19006
$ws.Range("M116").Value = -1471.2
$ws.Range("N116").Value = -25890
$ws.Range("H131").Value = 4917.294
$ws.Range("I131").Value = 674.2222
$ws.Range("J131").Value = 9690.75
$ws.Range("K131").Value = 2022.6666
$ws.Range("L131").Value = 29072.25
$ws.Range("M131").Value = 3017.3334
$ws.Range("N131").Value = -39152.25
$ws.Range("H132").Value = 1335.591
$ws.Range("I132").Value = 1078.7
$ws.Range("J132").Value = 3904.5
$ws.Range("K132").Value = 3236.1
$ws.Range("L132").Value = 11713.5
$ws.Range("M132").Value = -706.1000000000004
$ws.Range("N132").Value = -16773.5
$ws.Range("H138").Value = 2610
$ws.Range("I138").Value = 1485.9166
$ws.Range("J138").Value = 3025.0461
$ws.Range("K138").Value = 4457.7498
$ws.Range("L138").Value = 9075.138300000001
$ws.Range("M138").Value = 682.2502000000004
$ws.Range("N138").Value = -19355.1383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2622.6533
$ws.Range("I32").Value = 1995.7606
$ws.Range("K32").Value = 1995.7606
$ws.Range("M32").Value = -1708.7606
$ws.Range("H138").Value = 74999.5
$ws.Range("J138").Value = 74999.5
$ws.Range("L138").Value = 74999.5
$ws.Range("N138").Value = -85279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6433.875
$ws.Range("I20").Value = 5361.4
$ws.Range("K20").Value = 5361.4
$ws.Range("M20").Value = -5114.4
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("H99").Value = 3149.6667
$ws.Range("I99").Value = 3066
$ws.Range("J99").Value = 3233.3333
$ws.Range("K99").Value = 3066
$ws.Range("L99").Value = 3233.3333
$ws.Range("M99").Value = -1568
$ws.Range("N99").Value = -6229.3333
$ws.Range("H105").Value = 69874.75
$ws.Range("I105").Value = 101500
$ws.Range("K105").Value = 101500
$ws.Range("M105").Value = -99753
$ws.Range("H134").Value = 2127
$ws.Range("I134").Value = 2127
$ws.Range("K134").Value = 6381
$ws.Range("M134").Value = -3846
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3671.7693
$ws.Range("I99").Value = 3490
$ws.Range("J99").Value = 3785.375
$ws.Range("K99").Value = 3490
$ws.Range("L99").Value = 3785.375
$ws.Range("M99").Value = -1992
$ws.Range("N99").Value = -6781.375
$ws.Range("H126").Value = 3671.7693
$ws.Range("I126").Value = 3490
$ws.Range("J126").Value = 3785.375
$ws.Range("K126").Value = 10470
$ws.Range("L126").Value = 11356.125
$ws.Range("M126").Value = -8000
$ws.Range("N126").Value = -16296.125
$ws.Range("H132").Value = 4002.4
$ws.Range("I132").Value = 2902.3333
$ws.Range("K132").Value = 8706.999899999999
$ws.Range("M132").Value = -6176.999899999999
$ws.Range("H134").Value = 2518.2964
$ws.Range("I134").Value = 1999.28
$ws.Range("J134").Value = 9006
$ws.Range("K134").Value = 5997.84
$ws.Range("L134").Value = 27018
$ws.Range("M134").Value = -3462.84
$ws.Range("N134").Value = -32088
$ws.Range("H141").Value = 251444.5
$ws.Range("J141").Value = 251444.5
$ws.Range("L141").Value = 251444.5
$ws.Range("N141").Value = -261804.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 2998.3333
$ws.Range("I16").Value = 2990
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 8970
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -8797
$ws.Range("N16").Value = -9346
$ws.Range("H37").Value = 205059.23
$ws.Range("J37").Value = 205059.23
$ws.Range("L37").Value = 615177.6900000001
$ws.Range("N37").Value = -615401.6900000001
$ws.Range("H69").Value = 9252.75
$ws.Range("J69").Value = 9252.75
$ws.Range("L69").Value = 27758.25
$ws.Range("N69").Value = -29380.25
$ws.Range("H72").Value = 9252.75
$ws.Range("J72").Value = 9252.75
$ws.Range("L72").Value = 83274.75
$ws.Range("N72").Value = -91386.75
$ws.Range("H131").Value = 9128846
$ws.Range("I131").Value = 41667800
$ws.Range("J131").Value = 6078319.5
$ws.Range("K131").Value = 125003400
$ws.Range("L131").Value = 18234958.5
$ws.Range("M131").Value = -124998360
$ws.Range("N131").Value = -18245038.5
$ws.Range("H139").Value = 3431.9429
$ws.Range("I139").Value = 2075.8096
$ws.Range("J139").Value = 5466.143
$ws.Range("K139").Value = 6227.4288
$ws.Range("L139").Value = 16398.429
$ws.Range("M139").Value = -1087.4288
$ws.Range("N139").Value = -26678.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4907.222
$ws.Range("I126").Value = 3612.5
$ws.Range("J126").Value = 5277.143
$ws.Range("K126").Value = 10837.5
$ws.Range("L126").Value = 15831.429
$ws.Range("M126").Value = -8367.5
$ws.Range("N126").Value = -20771.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3674.5
$ws.Range("I22").Value = 1084.15
$ws.Range("K22").Value = 1084.15
$ws.Range("M22").Value = -789.1500000000001
$ws.Range("H27").Value = 3674.5
$ws.Range("I27").Value = 1084.15
$ws.Range("K27").Value = 1084.15
$ws.Range("M27").Value = -977.1500000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5836.3335
$ws.Range("I132").Value = 3025.7778
$ws.Range("K132").Value = 9077.3334
$ws.Range("M132").Value = -6547.3334
